$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "展览";     Cell = "F5";  Value = 431 },
    @{ Sheet = "展览";     Cell = "F6";  Value = 807 },
    @{ Sheet = "展览";     Cell = "F12"; Value = 673 },
    @{ Sheet = "展览";     Cell = "F13"; Value = 184 },
    @{ Sheet = "展览";     Cell = "F26"; Value = 5245 },
    @{ Sheet = "展览";     Cell = "F32"; Value = 1079 },
    @{ Sheet = "展览";     Cell = "F34"; Value = 51 },

    @{ Sheet = "演出";     Cell = "F17"; Value = 981 },
    @{ Sheet = "演出";     Cell = "F20"; Value = 615 },
    @{ Sheet = "演出";     Cell = "F25"; Value = 274 },
    @{ Sheet = "演出";     Cell = "F26"; Value = 3897 },

    @{ Sheet = "本地生活"; Cell = "F5";  Value = 2437 },
    @{ Sheet = "本地生活"; Cell = "F6";  Value = 1028 },

    @{ Sheet = "全部类型"; Cell = "F4";  Value = 2437 },
    @{ Sheet = "全部类型"; Cell = "F6";  Value = 1028 },
    @{ Sheet = "全部类型"; Cell = "F11"; Value = 431 },
    @{ Sheet = "全部类型"; Cell = "F17"; Value = 673 },
    @{ Sheet = "全部类型"; Cell = "F30"; Value = 5245 },
    @{ Sheet = "全部类型"; Cell = "F45"; Value = 274 },
    @{ Sheet = "全部类型"; Cell = "F46"; Value = 1079 },
    @{ Sheet = "全部类型"; Cell = "F50"; Value = 51 }
)

foreach ($change in $changes) {
    $ws = $wb.Worksheets.Item($change.Sheet)
    $ws.Range($change.Cell).Value = $change.Value
}
